# [Outlook] (internet headers) Add snippets
# Adds 5 new rows to the "Snippets" table (rows 252-256), growing the
# table/worksheet from A1:E251 to A1:E256, and moves the grid selection
# to the last new cell (E256) to mirror the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 fresh rows right after the current last data row (251) so the
# new rows inherit the existing body formatting (style "1") instead of
# picking up a brand-new style.
$ws.Rows("252:256").Insert() | Out-Null

# Row 252: MessageRead.getAllInternetHeadersAsync
$ws.Range("A252").Value = "MessageRead"
$ws.Range("B252").Value = "getAllInternetHeadersAsync"
$ws.Range("C252").Value = 1
$ws.Range("D252").Value = "outlook-mime-headers-get-internet-headers-message-read"
$ws.Range("E252").Value = "run"

# Row 253: MessageCompose.internetHeaders
$ws.Range("A253").Value = "MessageCompose"
$ws.Range("B253").Value = "internetHeaders"
$ws.Range("D253").Value = "outlook-mime-headers-manage-custom-internet-headers-message-compose"
$ws.Range("E253").Value = "getSelectedCustomHeaders"

# Row 254: InternetHeaders.getAsync
$ws.Range("A254").Value = "InternetHeaders"
$ws.Range("B254").Value = "getAsync"
$ws.Range("C254").Value = 1
$ws.Range("D254").Value = "outlook-mime-headers-manage-custom-internet-headers-message-compose"
$ws.Range("E254").Value = "getSelectedCustomHeaders"

# Row 255: InternetHeaders.removeAsync
$ws.Range("A255").Value = "InternetHeaders"
$ws.Range("B255").Value = "removeAsync"
$ws.Range("C255").Value = 1
$ws.Range("D255").Value = "outlook-mime-headers-manage-custom-internet-headers-message-compose"
$ws.Range("E255").Value = "removeSelectedCustomHeaders"

# Row 256: InternetHeaders.setAsync
$ws.Range("A256").Value = "InternetHeaders"
$ws.Range("B256").Value = "setAsync"
$ws.Range("C256").Value = 1
$ws.Range("D256").Value = "outlook-mime-headers-manage-custom-internet-headers-message-compose"
$ws.Range("E256").Value = "setCustomHeaders"

# Grow the "Snippets" table/autofilter to cover the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E256")) | Out-Null

# Move the active selection to the last cell that was edited.
$ws.Range("E256").Select() | Out-Null
